$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape Id=40 ("2 - Metodologia" text box) ---
$shp1 = $s.Shapes.Item(12)
$tr1 = $shp1.TextFrame.TextRange

# Paragraph: "Diante disso..." -> reworded comparison sentence
$para = $tr1.Characters(342, 240)
$para.Text = "`tDiante disso"
$anchor = $tr1.Characters(354, 1)
$anchor.InsertAfter(", ") | Out-Null
$anchor = $tr1.Characters(356, 1)
$anchor.InsertAfter("é possível ") | Out-Null
$anchor = $tr1.Characters(367, 1)
$anchor.InsertAfter("comparar ") | Out-Null
$anchor = $tr1.Characters(376, 1)
$anchor.InsertAfter("as sintaxes e estratégias de programação das linguagens com níveis diferentes") | Out-Null
$anchor = $tr1.Characters(453, 1)
$anchor.InsertAfter(", desde as que mais se aproximam da forma que o computador “compreende” o código até ") | Out-Null
$anchor = $tr1.Characters(538, 1)
$anchor.InsertAfter("as") | Out-Null
$anchor = $tr1.Characters(540, 1)
$anchor.InsertAfter(" que se ") | Out-Null
$anchor = $tr1.Characters(548, 1)
$anchor.InsertAfter("parecem ") | Out-Null
$anchor = $tr1.Characters(556, 1)
$anchor.InsertAfter("mais com a escrita humana.") | Out-Null

# Paragraph: "Para fins comparativos..." -> expanded palindrome explanation
$para = $tr1.Characters(18, 322)
$para.Text = "`tPara fins comparativos entre as diversas linguagens de programação existentes, realizamos uma aplicação para determinar se um conjunto de caracteres é palíndromo ou não, ou seja, "
$anchor = $tr1.Characters(197, 1)
$anchor.InsertAfter("quando a ") | Out-Null
$anchor = $tr1.Characters(206, 1)
$anchor.InsertAfter("palavra é ") | Out-Null
$anchor = $tr1.Characters(216, 1)
$anchor.InsertAfter("lida") | Out-Null
$anchor = $tr1.Characters(220, 1)
$anchor.InsertAfter(" ") | Out-Null
$anchor = $tr1.Characters(221, 1)
$anchor.InsertAfter("da direita para esquerda ou da esquerda para direita ") | Out-Null
$anchor = $tr1.Characters(274, 1)
$anchor.InsertAfter("e continua com a mesma sequencia de caracteres ") | Out-Null
$anchor = $tr1.Characters(321, 1)
$anchor.InsertAfter(", significa que ") | Out-Null
$anchor = $tr1.Characters(337, 1)
$anchor.InsertAfter(" se  trata  de  um palíndromo.") | Out-Null

# --- Shape Id=27 ("3 - Resultados" text box) ---
$shp2 = $s.Shapes.Item(21)
$tr2 = $shp2.TextFrame.TextRange

# Paragraph: "Assim, ao almejar..." -> split out " abstratos " run
$para = $tr2.Characters(17, 349)
$para.Text = "`t Assim, ao almejar o mesmo objetivo de verificação do palíndromo com linguagens menos usadas, atualmente, até as mais usadas, percebemos que alguns códigos são relativamente simples de serem compreendidos, já outros são mais "
$anchor = $tr2.Characters(242, 1)
$anchor.InsertAfter(" abstratos ") | Out-Null
$anchor = $tr2.Characters(253, 1)
$anchor.InsertAfter("e difíceis de serem entendidos porém, mesmo com diferenças significativas, é possível chegar ao mesmo resultado. ") | Out-Null
